$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-save the "plan_create_datetime" timestamps for PLN-001 rows (2-10).
# These carry the same instant but Excel re-serializes the underlying
# float with slightly different precision on save.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = 45403.689928449072
}

# --- Convert the PLN-002 rows (11-19) plan_create_datetime from text to a
# real datetime serial, matching the format used in column B elsewhere.
$ws.Cells.Item(2, 2).Copy()
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = 45403.690673573059
}

# --- Append 11 blank formatted rows (20-30) below the data, matching the
# D/H/I column date/datetime styles used elsewhere in the sheet.
$ws.Cells.Item(2, 4).Copy()
$ws.Range("D20:D30").PasteSpecial(-4122)

$ws.Cells.Item(2, 2).Copy()
$ws.Range("H20:H30").PasteSpecial(-4122)

$ws.Cells.Item(2, 4).Copy()
$ws.Range("I20:I30").PasteSpecial(-4122)

$ws.Range("D20:D30,H20:H30,I20:I30").ClearContents()

# --- Selection / view state as left by the editor.
$ws.Range("A20:K30").Select()
$excel.ActiveWindow.ScrollRow = 4
